$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.68"
$ws.Range("E2").Value = "'0.70%"
$ws.Range("D3").Value = "'26.84"
$ws.Range("E3").Value = "'-2.70%"
$ws.Range("D4").Value = "'4.729"
$ws.Range("E4").Value = "'-9.58%"
$ws.Range("E5").Value = "'0.81%"
$ws.Range("D6").Value = "'6.655"
$ws.Range("E6").Value = "'-0.96%"
$ws.Range("D7").Value = "'0.8683"
$ws.Range("E7").Value = "'0.55%"
$ws.Range("D8").Value = "'0.9467"
$ws.Range("E8").Value = "'-1.97%"
$ws.Range("B9").Value = "'WazirX"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'-0.54%"
$ws.Range("B10").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.03749"
$ws.Range("E10").Value = "'8.02%"
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07119"
$ws.Range("E11").Value = "'-0.63%"
$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03162"
$ws.Range("E12").Value = "'-0.56%"
$ws.Range("B13").Value = "'BitMartToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09257"
$ws.Range("E13").Value = "'0.26%"
$ws.Range("B14").Value = "'BitForexToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001541"
$ws.Range("E14").Value = "'-0.50%"
$ws.Range("B15").Value = "'One"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006036"
$ws.Range("E15").Value = "'-0.23%"
$ws.Range("D16").Value = "'0.006091"
$ws.Range("E16").Value = "'4.83%"
$ws.Range("D17").Value = "'3.498"
$ws.Range("E17").Value = "'-0.07%"
$ws.Range("D18").Value = "'3.202"
$ws.Range("E18").Value = "'-0.49%"
$ws.Range("E19").Value = "'-0.10%"
$ws.Range("E20").Value = "'-1.10%"
$ws.Range("D21").Value = "'0.1276"
$ws.Range("E21").Value = "'-1.75%"
$ws.Range("D22").Value = "'3.812"
$ws.Range("E22").Value = "'7.67%"
$ws.Range("D23").Value = "'0.04226"
$ws.Range("E23").Value = "'1.55%"
$ws.Range("D25").Value = "'0.001221"
$ws.Range("E25").Value = "'-0.45%"
$ws.Range("D26").Value = "'0.004289"
$ws.Range("E26").Value = "'-10.66%"
$ws.Range("D27").Value = "'0.0001189"
$ws.Range("E27").Value = "'-0.87%"
$ws.Range("D28").Value = "'0.0001492"
$ws.Range("E28").Value = "'1.77%"
$ws.Range("D40").Value = "'0.03824"
$ws.Range("E40").Value = "'0.20%"
$ws.Range("D41").Value = "'0.006207"
$ws.Range("E41").Value = "'9.29%"
$ws.Range("D42").Value = "'0.1103"
$ws.Range("E42").Value = "'0.12%"
$ws.Range("D43").Value = "'0.002163"
$ws.Range("E43").Value = "'-5.93%"
$ws.Range("D44").Value = "'0.01113"
$ws.Range("E44").Value = "'4.40%"
$ws.Range("D45").Value = "'0.00005496"
$ws.Range("E45").Value = "'4.92%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.06%"
$ws.Range("D47").Value = "'0.08846"
$ws.Range("E47").Value = "'-11.54%"
$ws.Range("D48").Value = "'0.002440"
$ws.Range("E48").Value = "'14.61%"
$ws.Range("E49").Value = "'-0.06%"
$ws.Range("E50").Value = "'-0.06%"
